$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A120").Value = "em_ui_filter"
$ws.Range("A121").Value = "em_ui_add"
$ws.Range("C121").Value = "追加 "
$ws.Range("D121").Value = "添加"
$ws.Range("D120").Value = "最近对话过滤"
$ws.Range("C120").Value = "最近の会話フィルター"
